$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "molecule" -> "protein complex" throughout column C (and a couple of
# column H cells) of the templates sheet, per commit:
#   - use GO:0042611 'MHC protein complex'
#   - use "protein complex" instead of "molecule" throughout

$ws.Range("C21").Value = '$class protein complex'
$ws.Range("H21").Value = 'MHC protein complex'

$ws.Range("C22").Value = '$taxon-label $class protein complex'
$ws.Range("C23").Value = '$prefix-$locus protein complex'
$ws.Range("C24").Value = '$prefix-$locus protein complex'
$ws.Range("C25").Value = '$prefix-$chain-ii-locus protein complex'
$ws.Range("C26").Value = '$molecule protein complex'

$ws.Range("C27").Value = '$taxon-label MHC protein complex with haplotype'
$ws.Range("C28").Value = '$prefix protein complex with $haplotype haplotype'
$ws.Range("C29").Value = '$class protein complex with serotype'
$ws.Range("C30").Value = '$taxon-label $class protein complex with serotype'
$ws.Range("C31").Value = '$prefix protein complex with $serotype serotype'

$ws.Range("C32").Value = 'mutant $class protein complex'
$ws.Range("H32").Value = 'mutant MHC protein complex'
$ws.Range("C33").Value = 'mutant $taxon-label $class protein complex'

$ws.Range("O34").Value = '$molecule protein complex'

# Widen column C to fit the new, longer "protein complex" labels.
$ws.Columns.Item(3).ColumnWidth = 30.1

# Move the selection, as recorded in the saved view state.
$ws.Range("J19").Select()
